# Fix the "汽車" (car) property sheet: the header row (row 1) had
# accidentally been filled with row 2's data values instead of proper
# column headers, and the sheet was missing the common trailing
# property columns (property_category, category, date, legislator_name,
# legislator_id, source_file, index) that the other property sheets
# already have.
#
# This inserts 7 new columns (H:N) and rewrites row 1 (headers) and
# fills in row 2 (data) for both the existing and the newly added
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Make room for the 7 new trailing columns, shifting nothing else.
$ws.Range("H1:N1").EntireColumn.Insert()

# --- Row 1: proper column headers ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: data, including the newly added columns ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# "2012-04-30" must stay a literal text string (not get auto-converted
# to a date serial number the way a plain .Value assignment would).
# Enter it as a text-producing formula, then collapse it down to its
# static value in place, which leaves the cell's existing style alone.
$ws.Range("J2").Formula = "=""2012-04-30"""
$ws.Range("J2").Copy()
$ws.Range("J2").PasteSpecial(-4163)
$ws.Range("K2").Value = "魏明谷"
$ws.Range("L2").Value = 980
$ws.Range("M2").Value = "tmp386d1"
$ws.Range("N2").Value = 31
